$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The donor submitted_id test value changed from "FOOBAR" to "XY_DONOR_ABCD".
# Setting this also drops the now-unused "FOOBAR" shared-string entry and
# shifts the remaining shared-string indices, same as the target workbook.
$ws.Range("A2").Value = "XY_DONOR_ABCD"

# Column widths were made explicit for columns A, D and E.
$ws.Columns.Item(1).ColumnWidth = 18.1640625
$ws.Columns.Item(4).ColumnWidth = 15.83203125
$ws.Columns.Item(5).ColumnWidth = 15.33203125

# Selection moved from F11 to A2.
$ws.Range("A2").Select()
